$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Espinaca, Femacal de La Calera) needs to be
# inserted as the new row 418, pushing the existing row 418 and everything
# below it down by one row (old row 418 -> 419, ..., old row 494 -> 495).
$ws.Rows("418:418").Insert()

# Populate the newly inserted row 418 with the new record's data.
$ws.Range("A418").Value = 3
$ws.Range("B418").Value = "Femacal de La Calera"
$ws.Range("C418").Value = "Coquimbo"
$ws.Range("D418").Value = 44995
$ws.Range("E418").Value = 5
$ws.Range("F418").Value = 100112012
$ws.Range("G418").Value = "Espinaca"
$ws.Range("H418").Value = "Sin especificar"
$ws.Range("I418").Value = "Primera"
$ws.Range("J418").Value = 140
$ws.Range("K418").Value = 6000
$ws.Range("L418").Value = 6500
$ws.Range("M418").Value = 6232
$ws.Range("N418").Value = "$/docena de atados (3 kilos)"
$ws.Range("O418").Value = "Provincia de Quillota"
$ws.Range("P418").Value = 2077
$ws.Range("Q418").Value = 3
$ws.Range("R418").Value = "Hortaliza"
